$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Section_A")
$ws2 = $wb.Worksheets.Item("Section_B")

# Section_A (sheet1) updates
$ws1.Range("B2").Value = "Free"
$ws1.Range("C2").Value = "Free"
$ws1.Range("D2").Value = "Free"
$ws1.Range("E2").Value = "DS401"
$ws1.Range("F2").Value = "DS401"
$ws1.Range("B3").Value = "DS456"
$ws1.Range("C3").Value = "DS456"
$ws1.Range("D3").Value = "DS456"
$ws1.Range("E3").Value = "DS456"
$ws1.Range("F3").Value = "Free"
$ws1.Range("B5").Value = "Free"
$ws1.Range("C5").Value = "Free"
$ws1.Range("D5").Value = "DS401"
$ws1.Range("E5").Value = "Free"
$ws1.Range("F5").Value = "Free"
$ws1.Range("B6").Value = "DS401"
$ws1.Range("C6").Value = "Free"
$ws1.Range("D6").Value = "Free"
$ws1.Range("E6").Value = "Free"
$ws1.Range("F6").Value = "Free"
$ws1.Range("B7").Value = "EC456"
$ws1.Range("C7").Value = "EC456"
$ws1.Range("D7").Value = "EC456"
$ws1.Range("E7").Value = "EC456"
$ws1.Range("F7").Value = "Free"

# Section_B (sheet2) updates
$ws2.Range("B2").Value = "DS456"
$ws2.Range("C2").Value = "DS456"
$ws2.Range("D2").Value = "Free"
$ws2.Range("E2").Value = "Free"
$ws2.Range("F2").Value = "DS401"
$ws2.Range("B3").Value = "DS401"
$ws2.Range("C3").Value = "Free"
$ws2.Range("D3").Value = "Free"
$ws2.Range("E3").Value = "DS456"
$ws2.Range("F3").Value = "Free"
$ws2.Range("B5").Value = "Free"
$ws2.Range("C5").Value = "Free"
$ws2.Range("D5").Value = "DS456"
$ws2.Range("E5").Value = "Free"
$ws2.Range("F5").Value = "Free"
$ws2.Range("B6").Value = "Free"
$ws2.Range("C6").Value = "EC456"
$ws2.Range("D6").Value = "Free"
$ws2.Range("E6").Value = "EC456"
$ws2.Range("F6").Value = "Free"
$ws2.Range("B7").Value = "EC456"
$ws2.Range("C7").Value = "DS401"
$ws2.Range("D7").Value = "DS401"
$ws2.Range("E7").Value = "Free"
$ws2.Range("F7").Value = "EC456"
